$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "[R] XIA(eXplainable AI) 패키지 중 DALEX로 변수 중요도 뽑기(classification)"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/R-XIAeXplainable-AI-%ED%8C%A8%ED%82%A4%EC%A7%80-%EC%A4%91-DALEX%EB%A1%9C-%EB%B3%80%EC%88%98-%EC%A4%91%EC%9A%94%EB%8F%84-%EB%BD%91%EA%B8%B0classification"

$ws.Range("D9").Value = "머신러닝 개발자?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/machine-learning-soft-engineer/#utm_source=rss&utm_medium=rss&utm_campaign=machine-learning-soft-engineer"
